$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string that Excel would otherwise auto-convert to a
# number (e.g. "1", "2020") so it is stored as genuine text, without
# leaving behind an extra/unused number-format style in styles.xml.
# We do this by building the literal as a text formula in a scratch
# cell, copying its computed (text) result, and pasting values-only
# into the destination - then clearing the scratch cell again.
function Set-TextLiteral($addr, $text) {
    $scratch = $ws.Range("ZZ9999")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# H2: REPORTDATE - text (contains dashes/colons, stays text naturally)
$ws.Range("H2").Value = "2020-09-30 00:00:00"

# I2: BASIC_EPS - numeric
$ws.Range("I2").Value = 0.5600000000000001

# K2: TOTAL_OPERATE_INCOME - numeric
$ws.Range("K2").Value = 199310474.48

# L2: PARENT_NETPROFIT - numeric
$ws.Range("L2").Value = 34815240.66

# N2: YSTZ - numeric (was blank)
$ws.Range("N2").Value = 10.2555047095

# O2: SJLTZ - numeric (was blank)
$ws.Range("O2").Value = 88.0915834674

# P2: BPS - numeric (was blank)
$ws.Range("P2").Value = 2.408789764032

# Q2: MGJYXJJE - numeric (was blank)
$ws.Range("Q2").Value = -0.787090161613

# R2: XSMLL - numeric
$ws.Range("R2").Value = 32.9923158587

# AB2: ISNEW - text "1" (looks numeric -> needs the text-literal helper)
Set-TextLiteral "AB2" "1"

# AC2: QDATE - text (contains "Q", stays text naturally)
$ws.Range("AC2").Value = "2020Q3"

# AD2: DATATYPE - text (contains Chinese chars, stays text naturally)
$ws.Range("AD2").Value = "2020年 三季报"

# AE2: DATAYEAR - text "2020" (looks numeric -> needs the text-literal helper)
Set-TextLiteral "AE2" "2020"

# AG2: EITIME - text (contains dashes/colons, stays text naturally)
$ws.Range("AG2").Value = "2020-12-22 16:06:38"
